{"js": "// Add a reviewer comment anchored to the phrase \"reconocidos por su esfuerzo\"\n// inside the paragraph that starts with \"En contraste, el \\u201cestallido social\\u201d ...\".\n//\n// The comment text matches the one added by Crist\\u00e1n Ayala in the source\n// document: a question about whether the connotation of \"esfuerzo\" shifted\n// from resilience to the courage of demanding change.\n\nconst body = context.document.body;\n\n// Find the unique occurrence of the exact phrase that should be wrapped by\n// the comment range (there is only one match in this document).\nconst searchResults = body.search(\"reconocidos por su esfuerzo\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the target phrase \"reconocidos por su esfuerzo\".');\n}\n\nconst targetRange = searchResults.items[0];\n\nconst commentText =\n  \"\\u00bfPodemos constatar un cambio en la connotaci\\u00f3n del esfuerzo desde \" +\n  \"resiliencia \\u2014sobreponernos a la adversidad\\u2014 a la valent\\u00eda de exigir y demandar?\";\n\ntargetRange.insertComment(commentText);\n\nawait context.sync();\n", "ps1": "# Add a reviewer comment anchored to the phrase \"reconocidos por su esfuerzo\"\n# inside the paragraph that starts with \"En contraste, el \"estallido social\" ...\".\n#\n# The comment reproduces the one added by Cristian Ayala in the source\n# document, questioning whether the connotation of \"esfuerzo\" shifted from\n# resilience to the courage of demanding change.\n\n$word.UserName = \"Cristi\u00e1n Ayala\"\n$word.UserInitials = \"CA\"\n\n$d = $word.ActiveDocument\n\n# Locate the unique occurrence of the target phrase.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"reconocidos por su esfuerzo\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find the target phrase \"reconocidos por su esfuerzo\".'\n}\n\n$commentText = \"\u00bfPodemos constatar un cambio en la connotaci\u00f3n del esfuerzo desde resiliencia \u2014sobreponernos a la adversidad\u2014 a la valent\u00eda de exigir y demandar?\"\n\n$d.Comments.Add($rng, $commentText) | Out-Null\n"}
